$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format while writing values that look like
# plain numbers, so Excel does not auto-convert them (e.g. "249.14" -> 249.14).
# Column D already holds free-form text (e.g. "37.018.13"), so after writing we
# restore the default "Normal" style to avoid leaving a stray number-format ref.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '37.028.33'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '2.060.02'
$ws.Range('E3').Value = '  -2.11%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '249.14'
$ws.Range('E5').Value = '  -1.62%  '
$ws.Range('D6').Value = '0.669'
$ws.Range('E6').Value = '  +1.07%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '55.44'
$ws.Range('E8').Value = '  +12.31%  '
$ws.Range('D9').Value = '60.39'
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('E11').Value = '  +7.23%  '
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '15.02'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').Value = '2.360.70'
$ws.Range('E14').Value = '  -1.99%  '
$ws.Range('D15').Value = '0.812'
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('D17').Value = '2.059.92'
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('D18').Value = '36.944.29'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '73.92'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  +11.51%  '
$ws.Range('D21').Value = '14.21'
$ws.Range('E21').Value = '  +6.53%  '
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').Value = '237.33'
$ws.Range('E23').Value = '  -1.90%  '
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('E25').Value = '  -3.50%  '
$ws.Range('D26').Value = '171.32'
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('E27').Value = '  -4.01%  '
$ws.Range('D28').Value = '20.07'
$ws.Range('E28').Value = '  -5.13%  '
$ws.Range('E29').Value = '  -1.92%  '
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('E34').Value = '  +5.90%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '0.0877'
$ws.Range('E35').Value = '  -5.33%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('E37').Value = '  -5.92%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  +22.08%  '
$ws.Range('D41').Value = '18.11'
$ws.Range('E41').Value = '  +7.12%  '
$ws.Range('D42').Value = '4.66'
$ws.Range('E42').Value = '  +60.06%  '
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('E44').Value = '  -2.16%  '
$ws.Range('E45').Value = '  -1.81%  '
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').Value = '4.35'
$ws.Range('E47').Value = '  +10.40%  '
$ws.Range('E48').Value = '  +6.63%  '
$ws.Range('D49').Value = '1.299.01'
$ws.Range('E49').Value = '  -3.82%  '
$ws.Range('E50').Value = '  -0.38%  '
$ws.Range('E51').Value = '  -3.38%  '

$ws.Range("D2:D51").Style = "Normal"
